$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F15").Style = "60% - Accent6"
$ws.Range("F15").NumberFormat = "0%"
